# Fix dataset: update hate-speech comment counts and journal/social labels
# for rows 3-108 of the "commenti odio per notizia" sheet, and append five
# new rows (104-108) that complete the "Strage di Cutro" breakdown.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Alluvioni Emilia Romagna','FanPage','Facebook',14),
    @('Alluvioni Emilia Romagna','FanPage','Instagram',7),
    @('Alluvioni Emilia Romagna','FanPage','YouTube',3),
    @('Alluvioni Emilia Romagna','Il Corriere Della Sera','Facebook',4),
    @('Alluvioni Emilia Romagna','Il Corriere Della Sera','Instagram',3),
    @('Alluvioni Emilia Romagna','Il Corriere Della Sera','YouTube',5),
    @('Alluvioni Emilia Romagna','La Repubblica','Facebook',15),
    @('Alluvioni Emilia Romagna','La Repubblica','Instagram',2),
    @('Alluvioni Emilia Romagna','La Repubblica','YouTube',6),
    @('Attacco di Hamas a Israele','FanPage','Facebook',1),
    @('Attacco di Hamas a Israele','FanPage','Instagram',16),
    @('Attacco di Hamas a Israele','FanPage','YouTube',2),
    @('Attacco di Hamas a Israele','Il Corriere Della Sera','Facebook',10),
    @('Attacco di Hamas a Israele','Il Corriere Della Sera','Instagram',11),
    @('Attacco di Hamas a Israele','Il Corriere Della Sera','YouTube',3),
    @('Attacco di Hamas a Israele','La Repubblica','Facebook',17),
    @('Attacco di Hamas a Israele','La Repubblica','Instagram',11),
    @('Attacco di Hamas a Israele','La Repubblica','YouTube',11),
    @('Geolier e i fischi a Sanremo','FanPage','Facebook',1),
    @('Geolier e i fischi a Sanremo','FanPage','Instagram',12),
    @('Geolier e i fischi a Sanremo','FanPage','YouTube',6),
    @('Geolier e i fischi a Sanremo','Il Corriere Della Sera','Facebook',1),
    @('Geolier e i fischi a Sanremo','Il Corriere Della Sera','Instagram',8),
    @('Geolier e i fischi a Sanremo','Il Corriere Della Sera','YouTube',14),
    @('Geolier e i fischi a Sanremo','La Repubblica','Instagram',41),
    @('Geolier e i fischi a Sanremo','La Repubblica','YouTube',5),
    @('Il caso Chiara Ferragni e Balocco','FanPage','Facebook',26),
    @('Il caso Chiara Ferragni e Balocco','FanPage','Instagram',27),
    @('Il caso Chiara Ferragni e Balocco','FanPage','YouTube',14),
    @('Il caso Chiara Ferragni e Balocco','Il Corriere Della Sera','Facebook',3),
    @('Il caso Chiara Ferragni e Balocco','Il Corriere Della Sera','Instagram',16),
    @('Il caso Chiara Ferragni e Balocco','Il Corriere Della Sera','YouTube',11),
    @('Il caso Chiara Ferragni e Balocco','La Repubblica','Facebook',12),
    @('Il caso Chiara Ferragni e Balocco','La Repubblica','Instagram',42),
    @('Il caso Chiara Ferragni e Balocco','La Repubblica','YouTube',8),
    @('Incidente Youtubers','FanPage','Facebook',39),
    @('Incidente Youtubers','FanPage','Instagram',36),
    @('Incidente Youtubers','FanPage','YouTube',11),
    @('Incidente Youtubers','Il Corriere Della Sera','Facebook',8),
    @('Incidente Youtubers','Il Corriere Della Sera','Instagram',39),
    @('Incidente Youtubers','Il Corriere Della Sera','YouTube',7),
    @('Incidente Youtubers','La Repubblica','Facebook',27),
    @('Incidente Youtubers','La Repubblica','Instagram',22),
    @('Incidente Youtubers','La Repubblica','YouTube',4),
    @('Invasione Ucraina','FanPage','Facebook',16),
    @('Invasione Ucraina','FanPage','Instagram',23),
    @('Invasione Ucraina','FanPage','YouTube',11),
    @('Invasione Ucraina','Il Corriere Della Sera','Facebook',10),
    @('Invasione Ucraina','Il Corriere Della Sera','Instagram',10),
    @('Invasione Ucraina','Il Corriere Della Sera','YouTube',8),
    @('Invasione Ucraina','La Repubblica','Facebook',9),
    @('Invasione Ucraina','La Repubblica','Instagram',15),
    @('Invasione Ucraina','La Repubblica','YouTube',11),
    @('L''implosione del sottomarino Titan','FanPage','Facebook',4),
    @('L''implosione del sottomarino Titan','FanPage','Instagram',3),
    @('L''implosione del sottomarino Titan','FanPage','YouTube',3),
    @('L''implosione del sottomarino Titan','Il Corriere Della Sera','Facebook',2),
    @('L''implosione del sottomarino Titan','Il Corriere Della Sera','Instagram',4),
    @('L''implosione del sottomarino Titan','Il Corriere Della Sera','YouTube',1),
    @('L''implosione del sottomarino Titan','La Repubblica','Facebook',1),
    @('L''implosione del sottomarino Titan','La Repubblica','Instagram',10),
    @('L''implosione del sottomarino Titan','La Repubblica','YouTube',6),
    @('L''incoronazione di Re Carlo','FanPage','Facebook',1),
    @('L''incoronazione di Re Carlo','FanPage','Instagram',12),
    @('L''incoronazione di Re Carlo','FanPage','YouTube',11),
    @('L''incoronazione di Re Carlo','Il Corriere Della Sera','Facebook',13),
    @('L''incoronazione di Re Carlo','Il Corriere Della Sera','Instagram',18),
    @('L''incoronazione di Re Carlo','Il Corriere Della Sera','YouTube',10),
    @('L''incoronazione di Re Carlo','La Repubblica','Facebook',6),
    @('L''incoronazione di Re Carlo','La Repubblica','Instagram',17),
    @('L''incoronazione di Re Carlo','La Repubblica','YouTube',3),
    @('L''omicidio di Giulia Cecchettin','FanPage','Facebook',7),
    @('L''omicidio di Giulia Cecchettin','FanPage','Instagram',4),
    @('L''omicidio di Giulia Cecchettin','FanPage','YouTube',1),
    @('L''omicidio di Giulia Cecchettin','Il Corriere Della Sera','Facebook',6),
    @('L''omicidio di Giulia Cecchettin','Il Corriere Della Sera','Instagram',10),
    @('L''omicidio di Giulia Cecchettin','Il Corriere Della Sera','YouTube',7),
    @('L''omicidio di Giulia Cecchettin','La Repubblica','Facebook',16),
    @('L''omicidio di Giulia Cecchettin','La Repubblica','Instagram',16),
    @('L''omicidio di Giulia Cecchettin','La Repubblica','YouTube',6),
    @('La cattura di Matteo Messina Denaro','FanPage','Facebook',10),
    @('La cattura di Matteo Messina Denaro','FanPage','Instagram',7),
    @('La cattura di Matteo Messina Denaro','FanPage','YouTube',4),
    @('La cattura di Matteo Messina Denaro','Il Corriere Della Sera','Facebook',4),
    @('La cattura di Matteo Messina Denaro','Il Corriere Della Sera','Instagram',11),
    @('La cattura di Matteo Messina Denaro','Il Corriere Della Sera','YouTube',7),
    @('La cattura di Matteo Messina Denaro','La Repubblica','Facebook',11),
    @('La cattura di Matteo Messina Denaro','La Repubblica','Instagram',22),
    @('La cattura di Matteo Messina Denaro','La Repubblica','YouTube',5),
    @('La morte di Silvio Berlusconi','FanPage','Facebook',1),
    @('La morte di Silvio Berlusconi','FanPage','Instagram',17),
    @('La morte di Silvio Berlusconi','FanPage','YouTube',5),
    @('La morte di Silvio Berlusconi','Il Corriere Della Sera','Facebook',4),
    @('La morte di Silvio Berlusconi','Il Corriere Della Sera','Instagram',25),
    @('La morte di Silvio Berlusconi','Il Corriere Della Sera','YouTube',4),
    @('La morte di Silvio Berlusconi','La Repubblica','Instagram',17),
    @('La morte di Silvio Berlusconi','La Repubblica','YouTube',12),
    @('Strage di Cutro','FanPage','Facebook',9),
    @('Strage di Cutro','FanPage','Instagram',32),
    @('Strage di Cutro','FanPage','YouTube',13),
    @('Strage di Cutro','Il Corriere Della Sera','Facebook',28),
    @('Strage di Cutro','Il Corriere Della Sera','Instagram',15),
    @('Strage di Cutro','Il Corriere Della Sera','YouTube',10),
    @('Strage di Cutro','La Repubblica','Facebook',39),
    @('Strage di Cutro','La Repubblica','Instagram',13),
    @('Strage di Cutro','La Repubblica','YouTube',14)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 3
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
}
